# "El ultimo commit de la batalla final"
# Fills in the test-case tracker (tester name, browser, project under test,
# and the three test-case rows) and updates the project title/comision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Project header text (Comision updated) -----------------------------
$ws.Range("B1").Value = "Proyecto Final CODERHOUSE PYTHON (Comision: 47770)"

# --- Header / meta block (rows 1-4) --------------------------------------
# Order matters: it drives the order new shared strings are appended, so
# cells are written in the same sequence the original author used.
$ws.Range("B3").Value = "Lucas Gak"          # Escrito por
$ws.Range("E1").Value = "Opera"              # Navegador
$ws.Range("E2").Value = "103.0.4928.34"      # Version
$ws.Range("E3").Value = "SMART TRAINING (Academia de deportes)"  # Descripcion
$ws.Range("B4").Value = "Lucas Gak"          # Probado por

# Descripcion cell (E3:F3, merged) gets an underlined font.
$ws.Range("E3:F3").Font.Underline = $true

# --- Test case rows (7-9) --------------------------------------------------
# Fecha column uses the same date for all three cases.
$ws.Range("B7").Value = 45220
$ws.Range("B8").Value = 45220
$ws.Range("B9").Value = 45220

# Case 1: user registration
$ws.Range("C7").Value = "Registro de usuario"
$ws.Range("D7").Value = "Registra tu usuario y te pide logearte con el mismo"
$ws.Range("E7").Value = "Registra tu usuario y te pide logearte con el mismo"

# Case 2: navbar links
$ws.Range("C8").Value = "Funcion de los links de la navbar"
$ws.Range("E8").Value = "Te dirige a distintas secciones de la pagina"
$ws.Range("D8").Value = "Dirigirte a distintas secciones de la pagina"

# Case 3: logout
$ws.Range("C9").Value = "Logout"
$ws.Range("D9").Value = "Te desloguea de tu cuenta y te ofrece volver a loguear o volver al inicio"
$ws.Range("E9").Value = "Te desloguea de tu cuenta y te ofrece volver a loguear o volver al inicio"

# Row 9 grew taller to fit the wrapped text.
$ws.Rows(9).RowHeight = 30

# --- Column widths, widened to fit the new content -------------------------
$ws.Columns("C:C").ColumnWidth = 54.666666666666664
$ws.Columns("D:D").ColumnWidth = 57.666666666666664
$ws.Columns("E:E").ColumnWidth = 50.833333333333336

# --- Selection ends on E9, matching the last cell touched -------------------
$ws.Range("E9").Select()
